$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 95; all existing rows 95-114 shift down to 96-115.
$ws.Range("A95").EntireRow.Insert()

# Populate the newly inserted row 95 with the new weekly price-report entry.
$ws.Range("A95").Value = 7
$ws.Range("B95").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C95").Value = "Ñuble"
$ws.Range("D95").Value = 45015
$ws.Range("E95").Value = 16
$ws.Range("F95").Value = "Fruta"
$ws.Range("G95").Value = 100108
$ws.Range("H95").Value = "Tropicales y subtropicales"
$ws.Range("I95").Value = 100108002
$ws.Range("J95").Value = "Mango"
$ws.Range("K95").Value = "Sin especificar"
$ws.Range("L95").Value = "Primera"
$ws.Range("M95").Value = 40
$ws.Range("N95").Value = 8000
$ws.Range("O95").Value = 8000
$ws.Range("P95").Value = 8000
$ws.Range("Q95").Value = '$/bandeja 4 kilos'
$ws.Range("R95").Value = "Perú"
$ws.Range("S95").Value = 2000
$ws.Range("T95").Value = 4
